# Menambahkan Biaya Admin | Fix Selected Dropdown Bulan Pemnbayaran | Menambahkan Popup Transaksi Berhasil
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix selected dropdown "bulan pembayaran" (row 1, col G) ---
$ws.Range("G1").Value = 3

# Column E holds phone numbers with leading zeros that must stay TEXT.
# Force text formatting before writing so the leading zero survives,
# then clear the formatting again so the cells end up back on the
# sheet's default (unstyled) look, same as every other cell.
$ws.Range("E2:E4").NumberFormat = "@"

# --- Row 2: replace student #2 (Anwar Ramdhan) with Nadia Hertisa Isnaeni Putri ---
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "Nadia Hertisa Isnaeni Putri"
$ws.Range("C2").Value = "hertisanadia44@gmail.com "
$ws.Range("D2").Value = "Komplek Permata Kopo C-189"
$ws.Range("E2").Value = "0043171547"
$ws.Range("F2").Value = 192010523
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = "smkn1ktp@01"

# --- Row 3: replace student #3 (Gilang Saputra) with Arianti Apriani Sagita ---
$ws.Range("A3").Value = 13
$ws.Range("B3").Value = "Arianti Apriani Sagita"
$ws.Range("C3").Value = "ariantiaprianisagita@gmail.com"
$ws.Range("D3").Value = "Kp. Pasanggrahan Rt 02 Rw 06 Kec. Pasirjambu"
$ws.Range("E3").Value = "0023620702"
$ws.Range("F3").Value = 192010505
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = "smkn1ktp@01"

# --- Row 4: replace student #4 (Raqhin Kusmanadinata) with Ajeng Nurfadillah ---
$ws.Range("A4").Value = 14
$ws.Range("B4").Value = "Ajeng Nurfadillah"
$ws.Range("C4").Value = "ajengnurfadilah@gmail.com"
$ws.Range("D4").Value = "Jln. Raya Sayuran Rt 08 Rw 07"
$ws.Range("E4").Value = "0034169559"
$ws.Range("F4").Value = 192010501
$ws.Range("G4").Value = 2
$ws.Range("I4").Value = "smkn1ktp@01"

# Drop the text-number-format hint again now the values are locked in as text.
$ws.Range("E2:E4").ClearFormats()

# --- Row 5 (Davin Albar) is removed entirely ---
$ws.Rows.Item(5).Delete()
